# Updated method for formatting pairwise comparison outputs
#
# The correlation table's "M" (mean) and "SD" (standard deviation) columns
# (B and C, rows 2-10) previously held placeholder values (0.00/1.00-ish).
# They are replaced here with the actual computed descriptive statistics
# for each of the nine variables. The pairwise correlation coefficients in
# columns D:K are left untouched.
#
# Values are written as literal text (matching the workbook's existing
# t="s" shared-string cells) rather than numbers: writing directly via
# Range.Value on a numeric-looking string auto-converts the cell to a
# number (and bumps its style to a quoted-text variant), which would
# disturb the sheet's styling. Instead we stage each value as a text
# formula in an unused scratch area, copy it, and paste-special "values
# only" into the target cell - this preserves the original cell style and
# keeps the stored type as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$means = @("6.27", "6.03", "5.03", "6.26", "5.40", "5.33", "5.07", "5.57", "5.62")
$sds   = @("0.94", "0.73", "1.31", "0.82", "1.34", "1.15", "1.02", "1.14", "1.17")

# Stage the new text values (as formulas returning string literals) in a
# scratch range well outside the table (AA:AB) so the real data in E:K is
# never disturbed.
for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Range("AA$row").Formula = "=""" + $means[$i] + """"
    $ws.Range("AB$row").Formula = "=""" + $sds[$i] + """"
}

# Copy the staged text values and paste only the values (not the formula)
# into the M/SD columns, preserving the destination cells' existing style.
$ws.Range("AA2:AB10").Copy()
$ws.Range("B2:C10").PasteSpecial(-4163)

# Clean up the scratch area.
$ws.Range("AA2:AB10").Clear()
